# Re-applies the refreshed coinranking.com snapshot (price + 1h change)
# used by the "Updated cryptos list" GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.168.30'
$ws.Range("E2").Value = '  +2.69%  '
$ws.Range("D3").Value = '2.537.59'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''594.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.17%  '
$ws.Range("D6").Value = '''177.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.00%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +2.40%  '
$ws.Range("D9").Value = '2.537.79'
$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("E10").Value = '  +2.85%  '
$ws.Range("E11").Value = '  +2.87%  '
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("D14").Value = '''27.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.07%  '
$ws.Range("D15").Value = '2.997.65'
$ws.Range("E15").Value = '  +1.98%  '
$ws.Range("D16").Value = '''0.0000180'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.73%  '
$ws.Range("D17").Value = '68.025.08'
$ws.Range("E17").Value = '  +2.72%  '
$ws.Range("D18").Value = '2.525.68'
$ws.Range("E18").Value = '  +2.17%  '
$ws.Range("D19").Value = '''8.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.43%  '
$ws.Range("E20").Value = '  +2.59%  '
$ws.Range("D21").Value = '''365.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.74%  '
$ws.Range("D22").Value = '''4.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("E23").Value = '  +2.70%  '
$ws.Range("E24").Value = '  +1.79%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '''71.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.74%  '
$ws.Range("D27").Value = '''10.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.61%  '
$ws.Range("D28").Value = '2.664.85'
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  +3.34%  '
$ws.Range("D31").Value = '''550.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.46%  '
$ws.Range("D32").Value = '''8.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.76%  '
$ws.Range("D33").Value = '''1.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.85%  '
$ws.Range("E34").Value = '  +2.94%  '
$ws.Range("E35").Value = '  -0.38%  '
$ws.Range("D36").Value = '''1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").Value = '''1.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = '''156.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.59%  '
$ws.Range("D39").Value = '''18.89'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.06%  '
$ws.Range("D40").Value = '''18.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.89%  '
$ws.Range("E41").Value = '  +1.46%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").Value = '''5.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.99%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '''1.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.87%  '
$ws.Range("D44").Value = '''2.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.48%  '
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("E46").Value = '  +1.84%  '
$ws.Range("D47").Value = '''147.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("E48").Value = '  +2.23%  '
$ws.Range("D49").Value = '0.0₆0279'
$ws.Range("E49").Value = '  +3.77%  '
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("D51").Value = '''0.0757'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.87%  '
